$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prices in column D are plain text in the source data (e.g. "52.030.15"),
# so every write is prefixed with a leading apostrophe to force Excel to
# keep them as text instead of auto-coercing them into numbers.

# --- Rows 34 & 35 swap places (Hedera <-> RenderToken) plus value updates ---
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'5.97"
$ws.Range("E34").Value = "  +5.40%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0934"
$ws.Range("E35").Value = "  +10.07%  "

# --- Price (D) and Volume(1h) (E) updates for remaining rows ---
$ws.Range("D2").Value = "'52.031.99"
$ws.Range("E2").Value = "  +1.13%  "

$ws.Range("D3").Value = "'2.883.63"
$ws.Range("E3").Value = "  +3.68%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'351.81"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").Value = "'111.46"
$ws.Range("E6").Value = "  +3.08%  "

$ws.Range("E7").Value = "  +1.53%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "'39.97"
$ws.Range("E10").Value = "  +1.82%  "

$ws.Range("D11").Value = "'0.0862"
$ws.Range("E11").Value = "  +3.31%  "

$ws.Range("D12").Value = "'0.136"
$ws.Range("E12").Value = "  +0.45%  "

$ws.Range("D13").Value = "'20.05"
$ws.Range("E13").Value = "  +0.87%  "

$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("D15").Value = "'3.333.23"
$ws.Range("E15").Value = "  +3.61%  "

$ws.Range("D16").Value = "'0.992"
$ws.Range("E16").Value = "  +7.15%  "

$ws.Range("D17").Value = "'2.876.68"
$ws.Range("E17").Value = "  +2.91%  "

$ws.Range("D18").Value = "'52.039.56"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").Value = "'7.71"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("E20").Value = "  +7.43%  "

$ws.Range("D21").Value = "'13.81"
$ws.Range("E21").Value = "  +3.65%  "

$ws.Range("D22").Value = "'0.0₃0982"
$ws.Range("E22").Value = "  +1.60%  "

$ws.Range("D23").Value = "'70.99"
$ws.Range("E23").Value = "  +0.67%  "

$ws.Range("D24").Value = "'270.08"
$ws.Range("E24").Value = "  +1.43%  "

$ws.Range("D25").Value = "'2.78"
$ws.Range("E25").Value = "  +0.64%  "

$ws.Range("D26").Value = "'26.33"
$ws.Range("E26").Value = "  +1.87%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  +0.20%  "

$ws.Range("E29").Value = "  +2.56%  "

$ws.Range("D30").Value = "'38.62"
$ws.Range("E30").Value = "  +4.25%  "

$ws.Range("E31").Value = "  +0.90%  "

$ws.Range("E32").Value = "  +2.50%  "

$ws.Range("D33").Value = "'53.15"
$ws.Range("E33").Value = "  +2.58%  "

$ws.Range("D36").Value = "'0.0459"
$ws.Range("E36").Value = "  +3.91%  "

$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("E38").Value = "  +6.38%  "

$ws.Range("D39").Value = "'18.61"
$ws.Range("E39").Value = "  +0.77%  "

$ws.Range("E40").Value = "  +3.76%  "

$ws.Range("D43").Value = "'22.48"
$ws.Range("E43").Value = "  +2.62%  "

$ws.Range("D44").Value = "'121.92"
$ws.Range("E44").Value = "  +1.47%  "

$ws.Range("E45").Value = "  +1.14%  "

$ws.Range("D46").Value = "'3.60"
$ws.Range("E46").Value = "  +7.17%  "

$ws.Range("D47").Value = "'2.193.19"
$ws.Range("E47").Value = "  +3.00%  "

$ws.Range("D48").Value = "'2.50"
$ws.Range("E48").Value = "  +7.42%  "

$ws.Range("D49").Value = "'0.268"
$ws.Range("E49").Value = "  +18.81%  "

$ws.Range("D50").Value = "'0.948"
$ws.Range("E50").Value = "  +5.95%  "

$ws.Range("D51").Value = "'5.50"
$ws.Range("E51").Value = "  +0.62%  "

Write-Host "Applied cryptos update"
